$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

$ws.Range("A38").Value = 41220
$ws.Range("B38").Value = 2.25
$ws.Range("C38").Value = 0.75
$ws.Range("D38").Value = "Creation of installer, test case tc07 put to operation"

$ws.Range("D38").Select()
